# Auto-generated edit script for _UI Language.xlsx
# Adds a 'Domain' sheet (translations for get_outfit_suggestion_use_case.dart)
# and fills in new profile/closet-overview/stats rows on the renamed 'Screen' sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Rename Sheet1 -> Screen ---
$ws1.Name = "Screen"

# --- Add the new Domain sheet right after Screen ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Domain"

# --- Extend Screen with newly-formatted rows (125-134), matching the existing
#     row look (border + wrap + vertical-center, row height 15.75) BEFORE writing
#     any values into them, since brand-new rows start out unformatted ---
$ws1.Range("A124:B124").Copy()
$ws1.Range("A125:B134").PasteSpecial(-4122)
for ($r = 125; $r -le 134; $r++) {
    $ws1.Rows($r).RowHeight = 15.75
}
$excel.CutCopyMode = 0

# --- Fill in the newly-populated rows on Screen (118-129) ---
$ws1.Range("A118").Value = "Unnamed"
$ws1.Range("B118").Value = "Chưa có tên"
$ws1.Range("A119").Value = "Edit profile"
$ws1.Range("B119").Value = "Chỉnh sửa thông tin"
$ws1.Range("A120").Value = "Profile"
$ws1.Range("B120").Value = "Trang cá nhân"
$ws1.Range("A121").Value = "Refresh page"
$ws1.Range("B121").Value = "Tải lại trang"
$ws1.Range("A122").Value = "Category"
$ws1.Range("B122").Value = "Danh mục"
$ws1.Range("A123").Value = "Color"
$ws1.Range("B123").Value = "Màu sắc"
$ws1.Range("A124").Value = "Season"
$ws1.Range("B124").Value = "Mùa"
$ws1.Range("A125").Value = "Occasion"
$ws1.Range("B125").Value = "Mục đích"
$ws1.Range("A126").Value = "Settings"
$ws1.Range("B126").Value = "Cài đặt"
$ws1.Range("A127").Value = "Closets overview"
$ws1.Range("B127").Value = "Tổng quan Tủ đồ"
$ws1.Range("A128").Value = "Statistics"
$ws1.Range("B128").Value = "Thống kê"
$ws1.Range("A129").Value = "No data for statistics"
$ws1.Range("B129").Value = "Chưa có dữ liệu để thống kê."

# --- Populate the Domain sheet ---
$ws2.Range("A1:B1").Merge()
$ws2.Range("A1").Value = "get_outfit_suggestion_use_case.dart"
$ws2.Range("B1").Value = ""
$ws2.Range("A2").Value = "Get weather by saved coordinates:"
$ws2.Range("B2").Value = "Lấy thời tiết theo tọa độ đã lưu:"
$ws2.Range("A3").Value = "Manual location data missing, reverting to default."
$ws2.Range("B3").Value = "Dữ liệu vị trí thủ công bị thiếu, quay về mặc định."
$ws2.Range("A4").Value = "Getting weather by auto-detecting location…"
$ws2.Range("B4").Value = "Lấy thời tiết theo vị trí tự động…"
$ws2.Range("A5").Value = "Location services are disabled, reverting to default."
$ws2.Range("B5").Value = "Dịch vụ vị trí đang tắt, quay về mặc định."
$ws2.Range("A6").Value = "Location services are disabled."
$ws2.Range("B6").Value = "Dịch vụ vị trí đang tắt."
$ws2.Range("A7").Value = "Location permission are denied, reverting to default."
$ws2.Range("B7").Value = "Không có quyền truy cập vị trí, quay về mặc định."
$ws2.Range("A8").Value = "Location permission are denied."
$ws2.Range("B8").Value = "Không có quyền truy cập vị trí."
$ws2.Range("A9").Value = "Failed to load weather for suggestions, using default."
$ws2.Range("B9").Value = "Lỗi khi lấy dữ liệu thời tiết cho gợi ý, sử dụng mặc định."
$ws2.Range("A10").Value = "Please add items to your closet to get suggestions."
$ws2.Range("B10").Value = "Vui lòng thêm đồ vào tủ để nhận gợi ý."

# --- Selections: Screen shows A1:B12 selected (not the active tab);
#     Domain is the active tab with B11 selected ---
$ws1.Range("A1:B12").Select()
$ws2.Select()
$ws2.Range("B11").Select()
